# Apply updated sensitivity values + column width tweaks to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column width adjustments (A: 7 -> 6.77734375, B: 20.33203125 -> 20.5546875, C: 6.88671875 -> 6.77734375) ---
$ws.Columns.Item(1).ColumnWidth = 6.77734375
$ws.Columns.Item(2).ColumnWidth = 20.5546875
$ws.Columns.Item(3).ColumnWidth = 6.77734375


# --- New values for row 2 (Breeze), row 3 (Gale), row 4 (Storm), columns D:AH ---
$row2 = @(
    40859.420632815905, 44187.036943834777, 47817.650747848835, 49873.506377285288,
    51120.813567305828, 50096.722198659234, 50084.160847639621, 50011.964425545011,
    49234.579054031994, 49934.123185134966, 49383.486396487031, 49416.487444799117,
    48569.843368840295, 49119.426303352782, 48386.631977142213, 48122.45031241879,
    48572.191940919001, 48902.63679647744,  50056.42591186364,  51544.62231777572,
    53295.086497048775, 55106.412265273699, 57035.834493365299, 59195.569177382757,
    61603.457081555141, 64278.985233287152, 67243.385076759383, 70519.710151723702,
    73793.579028871856, 77510.96537933698,  78516.39435377231
)

$row3 = @(
    40862.01489470368,  44201.626775360979, 47843.88446701513,  49914.041849709523,
    51179.050835423972, 50145.970265619981, 50146.052203292857, 50090.447243408475,
    49528.471662976452, 50544.039016916133, 50801.781098555133, 51400.160720835585,
    51082.369616023869, 52288.820890633375, 52085.823076441608, 52390.864019469431,
    53524.151723899937, 54516.484372305131, 56499.742904772596, 58907.114994096628,
    61654.775857585031, 64536.664484045694, 67573.856683828824, 70910.046134453631,
    74562.07264932667,  78547.046915107465, 82881.682724827755, 87581.968508913822,
    91517.894866016242, 96155.712867188748, 94937.812154554253
)

$row4 = @(
    40866.316920224788, 44236.36014906122,  47906.287099993322, 50010.773619032028,
    51317.74751145659,  50263.285464175278, 50293.446915616048, 50277.471133678104,
    49919.32687037759,  51229.6988342185,   52148.417755495495, 53279.13168361954,
    53484.277576406559, 55344.948066308381, 55708.095418743149, 56632.911874979945,
    58513.090819570898, 60257.555659773076, 63181.743818835806, 66652.923661231558,
    70583.141297985392, 74768.218495457375, 79198.225511313984, 84056.189404497884,
    89366.054495107528, 95150.494274590921, 101429.7771915824,  108220.2728036975,
    113110.91633989454, 119135.88353016907, 114304.84327530439
)

$startCol = 4  # column D
for ($i = 0; $i -lt $row2.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(2, $col).Value = $row2[$i]
    $ws.Cells.Item(3, $col).Value = $row3[$i]
    $ws.Cells.Item(4, $col).Value = $row4[$i]
}
